$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "status"

$ws.Range("A2").Value = "GHsarcev7PR0pGK3"
$ws.Range("B2").Value = "200"

$ws.Range("A3").Value = "zwWNZF9j99XJPayu"
$ws.Range("B3").Value = "200"

$ws.Range("A4").Value = "JgClaIA3srPc11g3"
$ws.Range("B4").Value = "200"
